# adding data + spatial comparison for contUS
#
# The NOAA_NCDC row (previously row 6, with no TempCovStart/TempCovEnd
# values) is moved up to row 3, pushing ERA5, GLDAS and GRIDMET down by
# one row each. GRIDMET's SpatRes value also changes from "Stations" to
# the newly-measured "0.04° x 0.04°" grid resolution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: NOAA_NCDC (moved up from row 6) ---
$ws.Range("A3").Value = "NOAA_NCDC"
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "Daily"
$ws.Range("E3").Value = "Global"
$ws.Range("F3").Value = "Stations"
$ws.Range("G3").Value = "T"
$ws.Range("H3").Value = "F"
$ws.Range("I3").Value = "F"
$ws.Range("J3").Value = "F"
$ws.Range("K3").Value = "F"
$ws.Range("L3").Value = "T"
$ws.Range("M3").Value = "F"
$ws.Range("N3").Value = "F"
$ws.Range("O3").Value = "T"

# --- Row 4: ERA5 (was row 3) ---
$ws.Range("A4").Value = "ERA5"
$ws.Range("B4").Value = 1981
$ws.Range("C4").Value = 2021
$ws.Range("D4").Value = "Hourly"
$ws.Range("E4").Value = "Global"
$ws.Range("F4").Value = "0.1° x 0.1°"
$ws.Range("G4").Value = "T"
$ws.Range("H4").Value = "T"
$ws.Range("I4").Value = "T"
$ws.Range("J4").Value = "T"
$ws.Range("K4").Value = "T"
$ws.Range("L4").Value = "T"
$ws.Range("M4").Value = "F"
$ws.Range("N4").Value = "F"
$ws.Range("O4").Value = "T"

# --- Row 5: GLDAS (was row 4) ---
$ws.Range("A5").Value = "GLDAS"
$ws.Range("B5").Value = 2000
$ws.Range("C5").Value = 2021
$ws.Range("D5").Value = "3-hourly"
$ws.Range("E5").Value = "Global"
$ws.Range("F5").Value = "0.25° x 0.25°"
$ws.Range("G5").Value = "T"
$ws.Range("H5").Value = "T"
$ws.Range("I5").Value = "T"
$ws.Range("J5").Value = "T"
$ws.Range("K5").Value = "T"
$ws.Range("L5").Value = "T"
$ws.Range("M5").Value = "T"
$ws.Range("N5").Value = "T"
$ws.Range("O5").Value = "T"

# --- Row 6: GRIDMET (was row 5); SpatRes updated to new 0.04 deg grid ---
$ws.Range("A6").Value = "GRIDMET"
$ws.Range("B6").Value = 1979
$ws.Range("C6").Value = 2021
$ws.Range("D6").Value = "Daily"
$ws.Range("E6").Value = "US"
$ws.Range("F6").Value = "0.04° x 0.04°"
$ws.Range("G6").Value = "T"
$ws.Range("H6").Value = "F"
$ws.Range("I6").Value = "F"
$ws.Range("J6").Value = "T"
$ws.Range("K6").Value = "T"
$ws.Range("L6").Value = "T"
$ws.Range("M6").Value = "F"
$ws.Range("N6").Value = "F"
$ws.Range("O6").Value = "F"

# Update the saved selection/active cell, as recorded in the sheet view.
$ws.Range("D13").Select()
